$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.05606233333333333
$ws.Range("H2").Value = 0.168187
$ws.Range("M2").Value = 1.089882
$ws.Range("N2").Value = 3.269646
$ws.Range("O2").Value = 0.03774844717592688
$ws.Range("P2").Value = 0.03774844717592687
$ws.Range("Q2").Value = 0.061101327978
$ws.Range("R2").Value = 0.549911951802
$ws.Range("S2").Value = 0.03774844717592688
$ws.Range("T2").Value = 0.03774844717592687

# Row 3
$ws.Range("G3").Value = 0.05606233333333333
$ws.Range("H3").Value = 0.168187
$ws.Range("O3").Value = 0.8393391727152114
$ws.Range("P3").Value = 0.8393391727152113
$ws.Range("Q3").Value = 1.358591992879667
$ws.Range("R3").Value = 12.227327935917
$ws.Range("S3").Value = 0.8393391727152114
$ws.Range("T3").Value = 0.8393391727152113

# Row 4
$ws.Range("G4").Value = 0.05606233333333333
$ws.Range("H4").Value = 0.168187
$ws.Range("M4").Value = 2.816943666666667
$ws.Range("N4").Value = 8.450831000000001
$ws.Range("O4").Value = 0.0975658366673901
$ws.Range("P4").Value = 0.09756583666739009
$ws.Range("Q4").Value = 0.1579244348218889
$ws.Range("R4").Value = 1.421319913397
$ws.Range("S4").Value = 0.0975658366673901
$ws.Range("T4").Value = 0.09756583666739009

# Row 5
$ws.Range("G5").Value = 0.05606233333333333
$ws.Range("H5").Value = 0.168187
$ws.Range("M5").Value = 0.7318113333333334
$ws.Range("N5").Value = 2.195434
$ws.Range("O5").Value = 0.0253465434414716
$ws.Range("P5").Value = 0.0253465434414716
$ws.Range("Q5").Value = 0.04102705090644444
$ws.Range("R5").Value = 0.369243458158
$ws.Range("S5").Value = 0.0253465434414716
$ws.Range("T5").Value = 0.0253465434414716
